# Fixes having a double from-path in Excel (#15)
#
# Adds two duplicate rows that both redirect "/duplicate-from-path" to
# "/resolve", reproducing the "double from-path" scenario referenced in the
# commit message, then widens column A so the new (longer) path values are
# fully visible and moves the active selection past the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter column B before column A so the shared-string table picks up
# "/resolve" (index 7) ahead of "/duplicate-from-path" (index 8), matching
# the order the strings were authored in.
$ws.Range("B9").Value = "/resolve"
$ws.Range("A9").Value = "/duplicate-from-path"

$ws.Range("B10").Value = "/resolve"
$ws.Range("A10").Value = "/duplicate-from-path"

# Widen column A (best-fit) so the newly added (longer) path values are
# fully visible, matching Excel's "AutoFit Column Width" result.
$ws.Columns("A").ColumnWidth = 17.8

# Move the active cell/selection below the newly entered rows.
$ws.Range("B11").Select()
